$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Testing/testrapport rows (23-35 table) ---------------------------------
# Rows 31 & 32: status moves from "Ej påbörjat" to "Klart" and gains a
# "Verklig tid" (actual time) figure.
$ws.Range("C31").Value = "Klart"
$ws.Range("E31").Value = 2

$ws.Range("C32").Value = "Klart"
$ws.Range("E32").Value = 2

$ws.Range("C33").Value = "Klart"
$ws.Range("E33").Value = 4

$ws.Range("C34").Value = "Klart"
$ws.Range("E34").Value = 2

# Set B35 ("Tankar och funderingar") before the two genuinely-new task
# names below, so the shared-string table keeps the same relative order
# as the source edit (new string right after the existing long comment).
$ws.Range("B35").Value = "Tankar och funderingar"
$ws.Range("C35").Value = "Klart"
$ws.Range("E35").Value = 6

# Task names updated to reflect actual testing work done this iteration.
$ws.Range("B33").Value = "Testspecifikation och testfall"
$ws.Range("B34").Value = "Testning och testrapport"

# --- Column widths (D & E) ---------------------------------------------------
# Target raw widths are 11.28515625 / 10.42578125, but this host's
# ColumnWidth setter only resolves to 1/6-character steps, so feed it the
# input that lands on the nearest reachable step (11.333333.. / 10.5).
$ws.Columns.Item(4).ColumnWidth = 10.5
$ws.Columns.Item(5).ColumnWidth = 9.666666666666666

# --- View / selection ---------------------------------------------------
# Source view scrolled so row 21 is the top visible row, with F33 as the
# active cell. Best-effort scroll (harmless if the host ignores it) plus
# the selection change, which is the part this host's sheetView persists.
$ws.Application.ActiveWindow.ScrollRow = 21
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("F33").Select()
